$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cell updates (values that Excel will not auto-convert to numbers)
$plainUpdates = @{
    "D2" = "69.341.84"
    "E2" = "  -1.32%  "
    "D3" = "3.527.44"
    "E3" = "  -2.67%  "
    "E4" = "  +0.29%  "
    "E5" = "  -0.29%  "
    "E6" = "  -3.87%  "
    "E7" = "  -2.77%  "
    "E8" = "  +0.05%  "
    "E9" = "  -1.69%  "
    "E10" = "  -4.15%  "
    "E11" = "  -4.63%  "
    "E12" = "  -6.69%  "
    "E13" = "  -3.24%  "
    "E14" = "  +13.38%  "
    "D15" = "4.093.79"
    "E15" = "  -2.52%  "
    "D16" = "69.426.51"
    "E16" = "  -1.45%  "
    "D17" = "3.553.70"
    "E17" = "  -1.85%  "
    "E18" = "  -6.20%  "
    "E19" = "  -3.88%  "
    "E20" = "  -0.89%  "
    "E21" = "  -3.17%  "
    "E22" = "  +3.50%  "
    "E23" = "  +3.15%  "
    "E24" = "  +1.87%  "
    "E25" = "  -5.65%  "
    "E26" = "  -4.02%  "
    "E27" = "  -6.07%  "
    "E28" = "  +0.65%  "
    "E29" = "  -2.40%  "
    "E30" = "  -7.07%  "
    "E31" = "  -5.44%  "
    "E32" = "  -4.20%  "
    "E33" = "  -5.91%  "
    "E34" = "  -2.11%  "
    "D35" = "3.787.91"
    "E35" = "  -3.57%  "
    "B36" = "PEPE"
    "C36" = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
    "D36" = "0.0₃0807"
    "E36" = "  -10.57%  "
    "B37" = "Dai"
    "C37" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "E37" = "  -0.07%  "
    "B38" = "Stacks"
    "C38" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "E38" = "  +3.32%  "
    "B39" = "Bittensor"
    "C39" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "E39" = "  -5.28%  "
    "B40" = "Fetch.AI"
    "C40" = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
    "E40" = "  -7.84%  "
    "E41" = "  -5.67%  "
    "E42" = "  +0.01%  "
    "E43" = "  -7.71%  "
    "E44" = "  +0.28%  "
    "E45" = "  +0.71%  "
    "E46" = "  -0.14%  "
    "E47" = "  -3.50%  "
    "E48" = "  -0.14%  "
    "E49" = "  -3.64%  "
    "B50" = "CoreDAO"
    "C50" = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
    "E50" = "  +68.57%  "
    "B51" = "Jupiter"
    "C51" = "https://coinranking.com/coin/qMgTxtv34+jupiter-jup"
    "E51" = "  +17.38%  "
}
foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Numeric-looking text cell updates: force text format so Excel keeps them as strings
# (matching the original workbook, which stores these as inline/shared strings, not numbers)
$numericTextUpdates = @{
    "D5" = "196.11"
    "D6" = "582.50"
    "D7" = "0.611"
    "D8" = "1.00"
    "D9" = "0.203"
    "D10" = "0.624"
    "D11" = "51.69"
    "D12" = "0.0000285"
    "D13" = "9.25"
    "D14" = "674.26"
    "D18" = "12.42"
    "D19" = "18.50"
    "D21" = "0.967"
    "D22" = "18.36"
    "D23" = "105.87"
    "D24" = "5.25"
    "D25" = "4.37"
    "D26" = "2.93"
    "D27" = "10.18"
    "D28" = "9.69"
    "D29" = "33.29"
    "D30" = "4.38"
    "D31" = "6.85"
    "D32" = "11.84"
    "D34" = "61.93"
    "D37" = "1.00"
    "D38" = "3.68"
    "D39" = "499.40"
    "D40" = "2.93"
    "D41" = "0.371"
    "D42" = "0.134"
    "D43" = "34.56"
    "D44" = "0.0458"
    "D45" = "2.90"
    "D47" = "0.136"
    "D49" = "8.33"
    "D50" = "2.78"
    "D51" = "1.75"
}
foreach ($ref in $numericTextUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextUpdates[$ref]
    $cell.Style = "Normal"
}
